$d = $word.ActiveDocument

# 1. Add an extra line break after "{{ formType }}"
$rng = $d.Content
$rng.Find.Execute("{{ formType }}", $true, $false, $false, $false, $false, $true, 1, $false, "{{ formType }}^l", 2) | Out-Null

# 2. After "{{ createdAt }}" add line breaks + modifiedAt/generatedAt placeholders
$rng2 = $d.Content
$rng2.Find.Execute("{{ createdAt }}", $true, $false, $false, $false, $false, $true, 1, $false, "{{ createdAt }}^l{{ modifiedAt }}^l{{ generatedAt }}", 2) | Out-Null

# 3. Adjust the first table's first two column widths (shrink col 1 by 2 twips,
#    grow col 2 by 2 twips) - mirrors the small re-layout recorded in the diff.
$t = $d.Tables(1)
$t.Columns(1).Width = 193.45
$t.Columns(2).Width = 13.7
